# page object for add emergency and documentation for login page
#
# Update the "VerifyAddEmergencyContact" test-data sheet:
#   - fix the contact name typo "tony stark" -> "tony starks"
#   - remove the now-redundant second data row ("bala")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VerifyAddEmergencyContact")

$ws.Range("C2").Value = "tony starks"

$ws.Rows.Item(3).Delete()

$ws.Range("C2").Select()
